# Update odds values in row 13 (Venados vs Tapatio) to reflect latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G13"  = 1.65
    "H13"  = 3.9
    "I13"  = 4.45
    "J13"  = 2.18
    "K13"  = 2.27
    "L13"  = 4.6
    "N13"  = 3.45
    "O13"  = 1.65
    "P13"  = 1.98
    "Q13"  = 2.55
    "R13"  = 1.39
    "U13"  = 1.7
    "V13"  = 1.93
    "W13"  = 7.7
    "X13"  = 8.25
    "Z13"  = 12.5
    "AA13" = 12.5
    "AB13" = 23
    "AD13" = 7.7
    "AF13" = 65
    "AH13" = 14
    "AI13" = 27
    "AJ13" = 14.5
    "AK13" = 75
    "AL13" = 40
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
